$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the team record (Wins/Losses/Ties) for each data row
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 77   # AD
    $ws.Cells.Item($r, 31).Value = 85   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
